$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (old row 3 -> becomes row 4).
$ws.Rows("3:3").Insert()

# Copy the formatting (font/border) of the row that was pushed down (now row 4)
# onto the newly inserted row 3, so it starts from the same base style.
$ws.Range("A4:H4").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 3 with its data:
# Facility | Global | Reference particle | Kinetic energy | 20 | MeV
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Facility"
$ws.Range("C3").Value = "Global"
$ws.Range("D3").Value = "Reference particle"
$ws.Range("E3").Value = "Kinetic energy"
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = "MeV"
# H3 stays blank.

# The new row's border differs from the row below it only by not having
# a bottom edge (left/right thin border carries through, top stays none).
for ($col = 1; $col -le 8; $col++) {
  $cell = $ws.Cells.Item(3, $col)
  $cell.Borders.Item(9).LineStyle = -4142
}

# Now update the row that was pushed down to row 4 with its new values
# (it keeps its existing style, only D/E/F/G change).
$ws.Range("D4").Value = "Vacuum chamber"
$ws.Range("E4").Value = "Mother volume radius"
$ws.Range("F4").Value = 0.75
$ws.Range("G4").Value = "m"
